# Updates odds values on the "Jogos da Semana" worksheet (row 3 and row 4)
# to reflect the latest FlashScore odds refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 updates ---
$ws.Range("G3").Value = 2.6
$ws.Range("H3").Value = 2.9
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 3.4
$ws.Range("L3").Value = 3.75
$ws.Range("O3").Value = 1.53
$ws.Range("P3").Value = 2.38
$ws.Range("Q3").Value = 2.03
$ws.Range("R3").Value = 1.83
$ws.Range("S3").Value = 2.7
$ws.Range("T3").Value = 1.44
$ws.Range("U3").Value = 4.1
$ws.Range("V3").Value = 1.22
$ws.Range("W3").Value = 5.5
$ws.Range("X3").Value = 1.14
$ws.Range("AC3").Value = 6.5
$ws.Range("AD3").Value = 11
$ws.Range("AE3").Value = 11
$ws.Range("AF3").Value = 26
$ws.Range("AG3").Value = 26
$ws.Range("AI3").Value = 6
$ws.Range("AN3").Value = 7
$ws.Range("AO3").Value = 13
$ws.Range("AP3").Value = 12
$ws.Range("AQ3").Value = 34
$ws.Range("AR3").Value = 29

# --- Row 4 updates ---
$ws.Range("V4").Value = 1.17
